$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts the old C:I data into D:J
# and (per the engine's behaviour observed) the new column inherits the
# wrap-text style (s="1") used by column B/old-C.
$ws.Columns("C").Insert()

# Give the new column C its own width (matches the "12-10m" value-box column);
# 15.25 is the character-unit value that this engine's pixel-quantized column
# model resolves closest to the author's stored width of 16.1640625.
$ws.Columns("C").ColumnWidth = 15.25

# Header for the newly inserted column
$ws.Range("C1").Value = "12-10m"

# Yes/No answers for the new "12-10m" column, row by row
$ws.Range("C2").Value  = "Yes"
$ws.Range("C3").Value  = "Yes"
$ws.Range("C4").Value  = "Yes"
$ws.Range("C5").Value  = "No"
$ws.Range("C6").Value  = "No"
$ws.Range("C7").Value  = "No"
$ws.Range("C8").Value  = "No"
$ws.Range("C9").Value  = "No"
$ws.Range("C10").Value = "No"
$ws.Range("C11").Value = "No"
$ws.Range("C12").Value = "No"
$ws.Range("C13").Value = "No"
$ws.Range("C14").Value = "No"
$ws.Range("C15").Value = "No"
$ws.Range("C16").Value = "No"
$ws.Range("C17").Value = "No"
$ws.Range("C18").Value = "No"
$ws.Range("C19").Value = "No"
$ws.Range("C20").Value = "Yes"
$ws.Range("C21").Value = "Yes"
$ws.Range("C22").Value = "Yes"
$ws.Range("C23").Value = "Yes"
$ws.Range("C24").Value = "Yes"
$ws.Range("C25").Value = "Yes"
$ws.Range("C26").Value = "Yes"
$ws.Range("C27").Value = "Yes"
$ws.Range("C28").Value = "No"
$ws.Range("C29").Value = "No"
$ws.Range("C30").Value = "No"
$ws.Range("C31").Value = "No"
$ws.Range("C32").Value = "No"
$ws.Range("C33").Value = "No"
$ws.Range("C34").Value = "No"

# Update the sheet selection to match the author's final cursor position
$ws.Range("C28:C34").Select()

# Match the author's slightly shorter window height
$excel.ActiveWindow.Height = 13880
